# ppt/theme/theme1.xml is the presentation's main theme (used by the Slide
# Master / all slides). It currently carries the "Integral" theme's
# "Red Violet" 12-slot colour scheme. The authored edit swaps it for the
# stock PowerPoint "Office Theme" colour scheme (the scheme that, before the
# edit, lived on the Notes Master's theme part). Font scheme and format
# scheme are already identical between the two theme parts, so only the
# colour-scheme slots need to change.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster

$master.ColorScheme.Colors(1).RGB  = 0          # dk1      000000
$master.ColorScheme.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$master.ColorScheme.Colors(3).RGB  = 6968388    # dk2      44546A
$master.ColorScheme.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$master.ColorScheme.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$master.ColorScheme.Colors(6).RGB  = 3243501    # accent2  ED7D31
$master.ColorScheme.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$master.ColorScheme.Colors(8).RGB  = 49407      # accent4  FFC000
$master.ColorScheme.Colors(9).RGB  = 12874308   # accent5  4472C4
$master.ColorScheme.Colors(10).RGB = 4697456    # accent6  70AD47
$master.ColorScheme.Colors(11).RGB = 12673797   # hlink    0563C1
$master.ColorScheme.Colors(12).RGB = 7491477    # folHlink 954F72
